$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.48), maa://25390 (96.18), maa://36681 (87.34)'
$ws.Range('L2').Value = '*maa://24633 (56.52), *maa://30515 (69.9), *maa://34787 (73.33), maa://39402 (91.8), ***maa://20792 (11.93), ***maa://29083 (27.78)'
$ws.Range('H3').Value = 'maa://21247 (98.56), *maa://22748 (60.0)'
$ws.Range('P3').Value = 'maa://21249 (94.47), maa://26254 (96.55)'
$ws.Range('X4').Value = '**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (89.47), ***maa://36683 (28.26)'
$ws.Range('D6').Value = 'maa://42407 (95.92)'
$ws.Range('AB6').Value = 'maa://22739 (91.38)'
$ws.Range('X7').Value = 'maa://22399 (95.36), *maa://22758 (74.29)'
$ws.Range('A8').Value = '更新日期：2025.02.16 13:17:39'
$ws.Range('X8').Value = 'maa://21411 (95.93)'
$ws.Range('D9').Value = 'maa://22765 (91.58), *maa://21915 (67.86)'
$ws.Range('T9').Value = '**maa://22866 (30.19), maa://26222 (98.04)'
$ws.Range('X9').Value = 'maa://26223 (97.89)'
$ws.Range('D10').Value = '***maa://25695 (18.72), ***maa://39951 (15.69), ***maa://34206 (20.0), ***maa://39243 (25.0), *maa://45271 (55.17)'
$ws.Range('P10').Value = 'maa://28977 (90.48), maa://36669 (86.36), *maa://23264 (61.82)'
$ws.Range('AF10').Value = '*maa://25021 (53.76), *maa://22733 (60.0), **maa://22761 (50.0)'
$ws.Range('D11').Value = 'maa://36707 (99.48)'
$ws.Range('L11').Value = 'maa://21287 (89.11)'
$ws.Range('X11').Value = 'maa://36713 (97.75)'
$ws.Range('D13').Value = 'maa://24999 (92.08), maa://36673 (93.24), maa://25001 (85.71)'
$ws.Range('X13').Value = 'maa://34957 (81.33), *maa://22768 (51.61)'
$ws.Range('L14').Value = 'maa://26245 (96.71), maa://21288 (96.3), maa://39841 (95.05), maa://36682 (97.44)'
$ws.Range('P14').Value = 'maa://23250 (98.74), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('D15').Value = '*maa://22743 (77.73), maa://22734 (84.17), *maa://30808 (64.18), **maa://36048 (45.0), maa://45058 (91.67)'
$ws.Range('P17').Value = 'maa://23890 (80.58), *maa://24940 (67.86)'
$ws.Range('D18').Value = 'maa://24570 (97.33)'
$ws.Range('AF18').Value = '*maa://24313 (59.28), **maa://29784 (46.43)'
$ws.Range('AB19').Value = '*maa://30709 (65.52), *maa://36668 (57.5)'
$ws.Range('D20').Value = 'maa://21432 (90.0), maa://25198 (93.58), *maa://20795 (50.77), maa://36680 (91.18)'
$ws.Range('L20').Value = 'maa://41331 (85.23)'
$ws.Range('P21').Value = 'maa://24381 (84.21)'
$ws.Range('X21').Value = 'maa://20110 (86.76), maa://34946 (90.91)'
$ws.Range('H22').Value = 'maa://25236 (95.6), **maa://21678 (48.94), **maa://22735 (42.86)'
$ws.Range('D23').Value = '***maa://28036 (28.77), *maa://41753 (58.82)'
$ws.Range('L23').Value = 'maa://39756 (95.56), maa://39875 (94.37)'
$ws.Range('P23').Value = 'maa://30587 (91.84), *maa://29748 (75.97), ***maa://29785 (16.18), *maa://37566 (76.32)'
$ws.Range('X24').Value = 'maa://29988 (84.58), maa://23504 (93.11), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (77.63), ***maa://22815 (23.08)'
$ws.Range('D25').Value = 'maa://29753 (95.13)'
$ws.Range('T25').Value = 'maa://20109 (92.09), maa://22545 (100.0), *maa://42915 (75.0)'
$ws.Range('X25').Value = '*maa://29890 (78.72)'
$ws.Range('AB25').Value = 'maa://31215 (87.72), maa://24516 (80.22), maa://26001 (87.5)'
$ws.Range('AB26').Value = 'maa://42235 (93.94)'
$ws.Range('D28').Value = 'maa://24465 (90.98), maa://25725 (83.72)'
$ws.Range('X28').Value = 'maa://39929 (90.62), maa://41749 (90.59), ***maa://39723 (13.89)'
$ws.Range('AF28').Value = 'maa://36660 (92.15), *maa://36701 (66.67)'
$ws.Range('AF29').Value = '*maa://24080 (68.85), maa://42865 (81.67), ***maa://34960 (8.33)'
$ws.Range('AB30').Value = 'maa://42979 (96.84), maa://45822 (100.0), *maa://45045 (80.0)'
$ws.Range('L31').Value = 'maa://35926 (93.43), maa://36258 (84.96), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.5), maa://36667 (97.56), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (95.76), maa://41108 (88.0), maa://41238 (97.09), maa://45523 (100.0)'
$ws.Range('T33').Value = 'maa://45558 (87.5)'
$ws.Range('T34').Value = 'maa://24526 (93.02)'
$ws.Range('L35').Value = 'maa://41296 (96.2)'
$ws.Range('K37').NumberFormat = '@'
$ws.Range('K37').Value = '3'
$ws.Range('L37').Value = 'maa://45718 (98.35), maa://45789 (100.0), maa://47069 (100.0)'
$ws.Range('AF38').Value = 'maa://36697 (86.26)'
$ws.Range('H39').Value = 'maa://36670 (88.89), maa://25199 (84.82), maa://30434 (91.36), ***maa://25036 (16.0), maa://45059 (81.25), *maa://44165 (66.67)'
$ws.Range('O39').NumberFormat = '@'
$ws.Range('O39').Value = '2'
$ws.Range('P39').Value = 'maa://24709 (91.5), maa://47093 (100.0)'
$ws.Range('T39').Value = 'maa://45788 (81.91), *maa://45790 (75.0), *maa://47079 (75.0)'
$ws.Range('P41').Value = '**maa://35616 (40.0), maa://43177 (90.91)'
$ws.Range('T44').Value = 'maa://39366 (88.57)'
$ws.Range('H46').Value = 'maa://35931 (92.09), maa://43901 (92.31)'
$ws.Range('H53').Value = 'maa://32534 (93.96), **maa://32434 (33.33)'
$ws.Range('H55').Value = 'maa://32532 (92.0)'
$ws.Range('H57').Value = 'maa://25176 (98.39)'
